# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.175.78"
$ws.Range("E2").Value = "  -4.15%  "

$ws.Range("D3").Value = "1.659.66"
$ws.Range("E3").Value = "  -2.69%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "'218.35"
$ws.Range("E5").Value = "  -2.39%  "

$ws.Range("D6").Value = "'0.5162"
$ws.Range("E6").Value = "  -2.68%  "

$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").Value = "'0.2590"
$ws.Range("E8").Value = "  -2.46%  "

$ws.Range("D9").Value = "'0.06456"
$ws.Range("E9").Value = "  -1.76%  "

$ws.Range("D10").Value = "'19.94"
$ws.Range("E10").Value = "  -3.93%  "

$ws.Range("D11").Value = "'0.07805"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("D12").Value = "1.657.95"
$ws.Range("E12").Value = "  -2.11%  "

$ws.Range("D13").Value = "'4.299"
$ws.Range("E13").Value = "  -4.77%  "

$ws.Range("D14").Value = "1.886.80"
$ws.Range("E14").Value = "  -2.79%  "

$ws.Range("E15").Value = "  -3.72%  "

$ws.Range("D16").Value = "0.0₅8061"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("E17").Value = "  -4.86%  "

$ws.Range("D18").Value = "26.201.06"
$ws.Range("E18").Value = "  -4.07%  "

$ws.Range("D19").Value = "'212.25"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "'4.420"
$ws.Range("E21").Value = "  -4.21%  "

$ws.Range("D23").Value = "'5.957"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("D25").Value = "'144.69"
$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("D26").Value = "'1.758"
$ws.Range("E26").Value = "  +2.63%  "

$ws.Range("D27").Value = "'0.1164"
$ws.Range("E27").Value = "  -3.13%  "

$ws.Range("D28").Value = "'6.978"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("D29").Value = "'15.86"
$ws.Range("E29").Value = "  -1.25%  "

$ws.Range("D30").Value = "'0.05285"
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("D31").Value = "'1.255"
$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("D32").Value = "'3.367"
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("D33").Value = "'3.219"
$ws.Range("E33").Value = "  -5.36%  "

$ws.Range("D34").Value = "'1.573"
$ws.Range("E34").Value = "  -4.01%  "

$ws.Range("D35").Value = "'2.762"
$ws.Range("E35").Value = "  -3.78%  "

$ws.Range("D36").Value = "'2.371"
$ws.Range("E36").Value = "  -1.85%  "

$ws.Range("D37").Value = "'0.9291"
$ws.Range("E37").Value = "  -1.63%  "

$ws.Range("D38").Value = "1.167.60"
$ws.Range("E38").Value = "  +12.46%  "

$ws.Range("D39").Value = "'0.5665"
$ws.Range("E39").Value = "  -2.02%  "

$ws.Range("D40").Value = "'0.01594"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D42").Value = "'0.8444"
$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").Value = "'5.695"
$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("D44").Value = "'100.45"
$ws.Range("E44").Value = "  -0.59%  "

$ws.Range("D45").Value = "1.796.94"

$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").Value = "'0.4536"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("D48").Value = "'55.90"
$ws.Range("E48").Value = "  -3.08%  "

$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").Value = "'7.895"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("D51").Value = "'0.05054"
$ws.Range("E51").Value = "  -3.35%  "
